$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7687332034111023
$ws.Range("B1").Value = 0.9861460328102112
$ws.Range("C1").Value = 1.485002040863037
$ws.Range("D1").Value = 2.14806866645813
$ws.Range("E1").Value = 1.586078286170959
